$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.723.25'
$ws.Range("E2").Value = '  +3.48%  '
$ws.Range("D3").Value = '2.446.06'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '576.51'
$ws.Range("E5").Value = '  +2.69%  '
$ws.Range("D6").Value = '145.80'
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("D9").Value = '2.444.26'
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("E10").Value = '  +2.33%  '
$ws.Range("D11").Value = '0.161'
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("E12").Value = '  +1.82%  '
$ws.Range("E13").Value = '  +3.07%  '
$ws.Range("D14").Value = '28.43'
$ws.Range("E14").Value = '  +9.34%  '
$ws.Range("E15").Value = '  +6.29%  '
$ws.Range("D16").Value = '2.889.13'
$ws.Range("E16").Value = '  +2.09%  '
$ws.Range("D17").Value = '62.656.39'
$ws.Range("E17").Value = '  +3.89%  '
$ws.Range("D18").Value = '2.430.57'
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").Value = '7.77'
$ws.Range("E19").Value = '  -3.49%  '
$ws.Range("D20").Value = '10.92'
$ws.Range("E20").Value = '  +3.01%  '
$ws.Range("D21").Value = '0.0₆0832'
$ws.Range("E21").Value = '  +204.73%  '
$ws.Range("D22").Value = '327.04'
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D24").Value = '2.02'
$ws.Range("E24").Value = '  +12.02%  '
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").Value = '65.54'
$ws.Range("E26").Value = '  +1.21%  '
$ws.Range("D27").Value = '643.25'
$ws.Range("E27").Value = '  +15.09%  '
$ws.Range("E28").Value = '  +16.21%  '
$ws.Range("D29").Value = '8.46'
$ws.Range("E29").Value = '  +6.18%  '
$ws.Range("D30").Value = '0.0₃0979'
$ws.Range("E30").Value = '  +5.30%  '
$ws.Range("D32").Value = '8.19'
$ws.Range("E32").Value = '  +1.89%  '
$ws.Range("E33").Value = '  +7.43%  '
$ws.Range("E34").Value = '  +3.91%  '
$ws.Range("E35").Value = '  +6.28%  '
$ws.Range("E36").Value = '  +2.82%  '
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '4.74'
$ws.Range("E38").Value = '  +3.39%  '
$ws.Range("D39").Value = '5.47'
$ws.Range("E39").Value = '  +7.37%  '
$ws.Range("D40").Value = '153.04'
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("D42").Value = '18.60'
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("E43").Value = '  +8.76%  '
$ws.Range("E44").Value = '  +5.90%  '
$ws.Range("D45").Value = '42.52'
$ws.Range("E45").Value = '  +1.89%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("E47").Value = '  +28.10%  '
$ws.Range("D48").Value = '144.56'
$ws.Range("E48").Value = '  +2.30%  '
$ws.Range("D49").Value = '3.60'
$ws.Range("E49").Value = '  +2.35%  '
$ws.Range("D50").Value = '20.62'
$ws.Range("E50").Value = '  +7.48%  '
$ws.Range("D51").Value = '0.606'
$ws.Range("E51").Value = '  +3.50%  '
